$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "60.870.63"
$ws.Range("E2").Value = "  -3.15%  "

$ws.Range("D3").Value = "2.918.43"
$ws.Range("E3").Value = "  -3.78%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "583.52"
$ws.Range("E5").Value = "  -1.46%  "

$ws.Range("D6").Value = "144.51"
$ws.Range("E6").Value = "  -5.65%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").Value = "2.917.02"
$ws.Range("E9").Value = "  -3.64%  "

$ws.Range("D10").Value = "6.81"
$ws.Range("E10").Value = "  +6.62%  "

$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  -3.33%  "

$ws.Range("D14").Value = "33.64"
$ws.Range("E14").Value = "  -5.28%  "

$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "3.402.40"
$ws.Range("E16").Value = "  -3.73%  "

$ws.Range("D17").Value = "60.835.63"
$ws.Range("E17").Value = "  -3.20%  "

$ws.Range("D18").Value = "6.74"
$ws.Range("E18").Value = "  -4.72%  "

$ws.Range("D19").Value = "2.919.51"
$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("D20").Value = "430.78"
$ws.Range("E20").Value = "  -4.69%  "

$ws.Range("D21").Value = "13.66"
$ws.Range("E21").Value = "  -4.39%  "

$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").Value = "7.13"
$ws.Range("E23").Value = "  -4.58%  "

$ws.Range("D24").Value = "80.60"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").Value = "10.87"
$ws.Range("E25").Value = "  -1.79%  "

$ws.Range("E26").Value = "  -4.37%  "

$ws.Range("E27").Value = "  -2.75%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "7.24"
$ws.Range("E30").Value = "  -4.25%  "

$ws.Range("E31").Value = "  -2.92%  "

$ws.Range("D32").Value = "2.17"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("D33").Value = "26.58"
$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("D35").Value = "0.0₃0874"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("E37").Value = "  -4.32%  "

$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  -4.37%  "

$ws.Range("D39").Value = "49.88"
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  -4.27%  "

$ws.Range("E42").Value = "  -4.20%  "

$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").Value = "41.38"
$ws.Range("E44").Value = "  -2.56%  "

$ws.Range("D45").Value = "377.99"
$ws.Range("E45").Value = "  -4.55%  "

$ws.Range("E46").Value = "  -2.85%  "

$ws.Range("D47").Value = "2.681.25"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").Value = "132.67"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D50").Value = "24.45"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("E51").Value = "  -1.55%  "
